$d = $word.ActiveDocument

# --- Change 1: paragraph containing "oferta subidas" / "my -offers" ---
# Fix "oferta subidas" -> "ofertas subidas" and drop the gramStart/gramEnd proofErr wrap
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("oferta subidas por el usuario", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "Target paragraph 1 not found" }
$rng1.Expand(4) | Out-Null
$rng1.InsertXML('<w:p w14:paraId="36950B7C" w14:textId="6B6EE6BF" w:rsidR="00AC6C2D" w:rsidRDefault="00F17642" w:rsidP="00F17642"><w:r><w:t xml:space="preserve">Cabe destacar </w:t></w:r><w:r w:rsidR="00AC6C2D"><w:t>que,</w:t></w:r><w:r><w:t xml:space="preserve"> si el usuario no tiene el rol de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>company</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00AC6C2D"><w:t xml:space="preserve"> o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00AC6C2D"><w:t>admin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00AC6C2D"><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> en la vista de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="003B7946"><w:t>homeOffer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> no aparecerá el botón para añadir una nueva oferta </w:t></w:r><w:r w:rsidR="003B7946"><w:t>o ver las ofertas subidas</w:t></w:r><w:r w:rsidR="00CC1084"><w:t>.</w:t></w:r><w:r w:rsidR="0031454E"><w:t xml:space="preserve"> La vista para ver las ofertas subidas se llama </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="0031454E"><w:t>my</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="0031454E"><w:t xml:space="preserve"> -</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="0031454E"><w:t>offers</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="0031454E"><w:t xml:space="preserve"> y solo muestra las </w:t></w:r><w:r w:rsidR="0031454E"><w:t>ofertas subidas</w:t></w:r><w:r w:rsidR="0031454E"><w:t xml:space="preserve"> por el usuario pero se ve cómo la vista de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="0031454E"><w:t>homeOffer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="0031454E"><w:t>.</w:t></w:r></w:p>
')

# --- Change 2: paragraph starting "Si el usuario es el propietario..." ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("i el usuario es el propietario de la oferta o bien administrador", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Target paragraph 2 not found" }
$rng2.Expand(4) | Out-Null
$rng2.InsertXML('<w:p w14:paraId="481D7C76" w14:textId="76FB9A64" w:rsidR="00AC6C2D" w:rsidRDefault="00AC6C2D" w:rsidP="00F17642"><w:r><w:t xml:space="preserve">Si el usuario es el propietario de la oferta o bien administrador, este podrá editar la oferta y acceder a la vista </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>offerEdit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> cómo se puede ver en la figura NUMERO3. Esta vista permite al propietario cambiar los campos de la oferta, pero no puede cambiar el documento asociado a ella. Esto se debe a que el fichero de la oferta es el elemento principal al que van asociados el resto de los campos, es decir sin documento, dichos campos no irían asociados a nada, por ello, si se desea cambiar el fichero de la oferta, esta debe crearse desde cero y subir una nueva oferta.</w:t></w:r><w:r w:rsidR="007D0585"><w:t xml:space="preserve"> La vista cuenta además con un botón para eliminar la oferta.</w:t></w:r></w:p>')

# --- Change 3: paragraph starting "A la diferencia" ---
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("A la diferencia", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) { throw "Target paragraph 3 not found" }
$rng3.Expand(4) | Out-Null
$rng3.InsertXML('<w:p w14:paraId="264F5882" w14:textId="5734914E" w:rsidR="00CC1084" w:rsidRDefault="00CC1084" w:rsidP="00CC1084"><w:r><w:lastRenderedPageBreak/><w:t>A la diferencia de lo que ocurre con los v</w:t></w:r><w:r w:rsidR="007D0585"><w:t>í</w:t></w:r><w:r><w:t xml:space="preserve">deos donde se guardan dos ficheros, aquí solo se guarda un fichero, por ello no hace falta crear subcarpetas para cada uno dentro de la carpeta </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>offers</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Para guardar el fichero se utiliza el siguiente formato dentro de la carpeta </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>offers</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>:</w:t></w:r></w:p>')

# --- Change 4: paragraph ending "centrarse solo en las que les interesan." gets a trailing run + a new paragraph ---
$rng4 = $d.Content
$found4 = $rng4.Find.Execute("En esta vista", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found4) { throw "Target paragraph 4 not found" }
$rng4.Expand(4) | Out-Null
$rng4.InsertXML('<w:p w14:paraId="252463D7" w14:textId="642C1DCC" w:rsidR="00F17642" w:rsidRDefault="00F17642" w:rsidP="00F17642"><w:r><w:t xml:space="preserve">En esta vista </w:t></w:r><w:r w:rsidR="007D29E2"><w:t>encuentran</w:t></w:r><w:r><w:t xml:space="preserve"> todas las ofertas unas debajo de otras y se podrá</w:t></w:r><w:r w:rsidR="007D29E2"><w:t>n</w:t></w:r><w:r><w:t xml:space="preserve"> realizar búsquedas por título gracias a la barra de búsqueda</w:t></w:r><w:r w:rsidR="007D29E2"><w:t xml:space="preserve">, la cual permite que los usuarios </w:t></w:r><w:r><w:t xml:space="preserve">puedan </w:t></w:r><w:r w:rsidR="007D29E2"><w:t>centrarse solo en las que les interesan.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t>También existe otra barra de búsqueda que permite buscar las ofertas según un valor mínimo de salario (es decir, todas aquellas ofertas que tengan ese valor o uno mayor para el campo salario serán mostradas al usuario). De esta forma, los usuarios podrán buscar las ofertas por su título o por su salario, siendo ambas búsquedas compatibles, ya que se pueden buscar ofertas de un determinado título y dentro de ese título filtrar para ver aquellas que ofrecen las condiciones de salario que decida buscar el usuario.</w:t></w:r></w:p>')

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
